$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10-11: I2C block (merged A10:A11)
$ws.Range("A10").Value = "I2C"
$ws.Range("B10").Value = "i2c_scl"
$ws.Range("C10").Value = "N1"
$ws.Range("D10").Value = "LVCMOS33"
$ws.Range("E10").Value = "OUT"
$ws.Range("F10").Value = "3.3V"

$ws.Range("B11").Value = "i2c_sda"
$ws.Range("C11").Value = "P2"
$ws.Range("D11").Value = "LVCMOS33"
$ws.Range("E11").Value = "INOUT"
$ws.Range("F11").Value = "3.3V"

$ws.Range("A10:A11").Merge() | Out-Null

# Rows 12-19: Logic block (merged A12:A19)
$ws.Range("A12").Value = "Logic"
$ws.Range("B12").Value = "dc_signal_in[0]"
$ws.Range("C12").Value = "T18"
$ws.Range("D12").Value = "LVCMOS33"
$ws.Range("E12").Value = "IN"
$ws.Range("F12").Value = "3.3V"

$ws.Range("B13").Value = "dc_signal_in[1]"
$ws.Range("C13").Value = "T17"
$ws.Range("D13").Value = "LVCMOS33"
$ws.Range("E13").Value = "IN"
$ws.Range("F13").Value = "3.3V"

$ws.Range("B14").Value = "dc_signal_in[2]"
$ws.Range("C14").Value = "N14"
$ws.Range("D14").Value = "LVCMOS33"
$ws.Range("E14").Value = "IN"
$ws.Range("F14").Value = "3.3V"

$ws.Range("B15").Value = "dc_signal_in[3]"
$ws.Range("C15").Value = "M14"
$ws.Range("D15").Value = "LVCMOS33"
$ws.Range("E15").Value = "IN"
$ws.Range("F15").Value = "3.3V"

$ws.Range("B16").Value = "dc_signal_in[4]"
$ws.Range("C16").Value = "N16"
$ws.Range("D16").Value = "LVCMOS33"
$ws.Range("E16").Value = "IN"
$ws.Range("F16").Value = "3.3V"

$ws.Range("B17").Value = "dc_signal_in[5]"
$ws.Range("C17").Value = "N15"
$ws.Range("D17").Value = "LVCMOS33"
$ws.Range("E17").Value = "IN"
$ws.Range("F17").Value = "3.3V"

$ws.Range("B18").Value = "dc_signal_in[6]"
$ws.Range("C18").Value = "M18"
$ws.Range("D18").Value = "LVCMOS33"
$ws.Range("E18").Value = "IN"
$ws.Range("F18").Value = "3.3V"

$ws.Range("B19").Value = "dc_signal_in[7]"
$ws.Range("C19").Value = "M16"
$ws.Range("D19").Value = "LVCMOS33"
$ws.Range("E19").Value = "IN"
$ws.Range("F19").Value = "3.3V"

$ws.Range("A12:A19").Merge() | Out-Null

# Match the style of column A (s="2") for the new A cells
$ws.Range("A10:A19").HorizontalAlignment = -4108
$ws.Range("A10:A19").VerticalAlignment = -4108

# Match the style of columns B/D/E/F (s="1", center horizontal) for the new cells
$ws.Range("B10:B19").HorizontalAlignment = -4108
$ws.Range("D10:D19").HorizontalAlignment = -4108
$ws.Range("E10:E19").HorizontalAlignment = -4108
$ws.Range("F10:F19").HorizontalAlignment = -4108

# Column C: rows 10-11 centered (matching the existing column style),
# but rows 12-19 keep the default (general) alignment as in the source.
$ws.Range("C10:C11").HorizontalAlignment = -4108

$ws.Range("Q18").Select() | Out-Null
